$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data: updated B/C/D values for existing rows 3-13, plus new rows
# 14-17 ("Speed meta analysis" section) for both worksheets.
# ---------------------------------------------------------------------------

function Set-Row($ws, [int]$r, [double]$b, [double]$c, [double]$d) {
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}

# -------------------- Sheet 1 ("treatment") -------------------------------
$ws1 = $wb.Worksheets.Item("treatment")

Set-Row $ws1 3  43.41501637732537  62.2155919053358   86.40885476689654
Set-Row $ws1 4  0.08305185177882789 0.3044734882949733 0.6670767227068665
Set-Row $ws1 5  0.2970451281813159 0.5687514478983436 0.8418518430936012

Set-Row $ws1 7  13.55589194560701  22.39109759868253  34.92806484539663
Set-Row $ws1 8  0.1569854918691213 0.6004271768547113 1.335795430571468
Set-Row $ws1 9  0.408987775409322  0.7998540327668482 1.193027748302593

Set-Row $ws1 11 9.115531553365228  10.88505355762126  12.88088304512077
Set-Row $ws1 12 0.01508900280650253 0.06432127698288349 0.148570837962344
Set-Row $ws1 13 0.1272758780208357 0.2627805159770388 0.3993766028508818

$ws1.Cells.Item(14, 1).Value = "Speed meta analysis"

$ws1.Cells.Item(15, 1).Value = "mean (km/day)"
Set-Row $ws1 15 3.187781843134013 3.714053319368254 4.301081131839839

$ws1.Cells.Item(16, 1).Value = "CoV² (RVAR)"
Set-Row $ws1 16 0.01361039429577891 0.05052841464647288 0.1112129518881235

$ws1.Cells.Item(17, 1).Value = "CoV  (RSTD)"
Set-Row $ws1 17 0.120301390926437 0.2317945023887313 0.3438847994087175

# -------------------- Sheet 2 ("control") ----------------------------------
$ws2 = $wb.Worksheets.Item("control")

Set-Row $ws2 3  0.8924377766345891 1.482719831771361 2.318365442301688
Set-Row $ws2 4  0.03580327817873998 0.2916148262838491 0.8097097943181886
Set-Row $ws2 5  0.2005982807328768 0.5724935662415689 0.9539607500205403

Set-Row $ws2 7  29.38287136681903  40.33554768163637  54.03942886010705
Set-Row $ws2 8  0.017706117717383  0.1202710342982784 0.3190116413560197
Set-Row $ws2 9  0.1401627180412515 0.3653014196879659 0.5949412792327863

Set-Row $ws2 11 7.854514005431656  11.5368521042755   16.31705476738479
Set-Row $ws2 12 0.01695684690716749 0.1584788789454684 0.4544507422679111
Set-Row $ws2 13 0.1387220680607462 0.4240906525570559 0.7181520649451167

$ws2.Cells.Item(14, 1).Value = "Speed meta analysis"

$ws2.Cells.Item(15, 1).Value = "mean (km/day)"
Set-Row $ws2 15 3.890621297764927 4.502646334156091 5.180288057555789

$ws2.Cells.Item(16, 1).Value = "CoV² (RVAR)"
Set-Row $ws2 16 0.002988245315773836 0.02368412071564556 0.06532771068280042

$ws2.Cells.Item(17, 1).Value = "CoV  (RSTD)"
Set-Row $ws2 17 0.05789683738708705 0.1629954008527747 0.2707044761453327
